# Apply test data adjustments for "private Krankenkasse - ohne Abkuerzung.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update Umlage U1 in Prozent (B4): 1.6 -> 1.9
$ws.Range("B4").Value = 1.9

# Update Umlage U2 in Prozent (B5): 0.44 -> 0.39
$ws.Range("B5").Value = 0.39

# Update Eintragungsdatum (B7): 15.12.2023 -> 01.01.2024 (stored as text, cell is text-formatted)
$ws.Range("B7").Value = "01.01.2024"

# Move the active selection to B2 (as reflected in the saved sheet view)
$ws.Range("B2").Select()
